# Update gh-pages output: remove the "2024-03-16 南宁·草莓动漫节" entry (event has
# passed / been superseded) from both the "展览" and "全部类型" sheets, shifting the
# remaining rows up and refreshing the "想去人数" (want-to-go) counts that bilibili
# reported for the still-upcoming events.

$wb = $excel.ActiveWorkbook

function Update-Sheet($ws) {
    # Row 2 (the 2024-03-16 草莓动漫节 entry) is gone; Excel shifts rows 3..N up to 2..N-1.
    $ws.Rows(2).Delete()

    # Re-sequence the leading index column (A) so it stays 1,2,3,... after the shift.
    $lastRow = $ws.Cells(1, 1).End(-4121).Row
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }

    # Refresh "想去人数" (column F) for the events whose counts changed in this scrape.
    for ($r = 2; $r -le $lastRow; $r++) {
        $link = $ws.Cells.Item($r, 8).Value()
        if ($link -eq "https://show.bilibili.com/platform/detail.html?id=81658") {
            $ws.Cells.Item($r, 6).Value = 935
        } elseif ($link -eq "https://show.bilibili.com/platform/detail.html?id=82416") {
            $ws.Cells.Item($r, 6).Value = 1771
        } elseif ($link -eq "https://show.bilibili.com/platform/detail.html?id=82241") {
            $ws.Cells.Item($r, 6).Value = 403
        }
    }
}

foreach ($ws in $wb.Worksheets) {
    if (($ws.Name -eq "展览") -or ($ws.Name -eq "全部类型")) {
        Update-Sheet $ws
    }
}
